$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '36.501.57'
$ws.Range('E2').Value = '  +1.80%  '

# Row 3
$ws.Range('D3').Value = '1.949.29'
$ws.Range('E3').Value = '  -0.32%  '

# Row 4
$ws.Range('E4').Value = '  -0.31%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '243.75'
$ws.Range('E5').Value = '  +1.00%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.618'
$ws.Range('E6').Value = '  -1.02%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.48'
$ws.Range('E7').Value = '  -5.95%  '

# Row 8
$ws.Range('E8').Value = '  -0.23%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.368'
$ws.Range('E9').Value = '  -0.20%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '55.74'
$ws.Range('E10').Value = '  -0.50%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0840'
$ws.Range('E11').Value = '  +5.29%  '

# Row 12
$ws.Range('E12').Value = '  +1.13%  '

# Row 13
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.76'
$ws.Range('E13').Value = '  -0.95%  '

# Row 14
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.828'
$ws.Range('E14').Value = '  -3.11%  '

# Row 15
$ws.Range('D15').Value = '2.234.35'
$ws.Range('E15').Value = '  -1.00%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '13.69'
$ws.Range('E16').Value = '  -2.26%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.27'
$ws.Range('E17').Value = '  -2.48%  '

# Row 18
$ws.Range('D18').Value = '1.954.06'
$ws.Range('E18').Value = '  -1.11%  '

# Row 19
$ws.Range('D19').Value = '36.389.11'
$ws.Range('E19').Value = '  +1.56%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.90'
$ws.Range('E20').Value = '  -1.36%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0867'
$ws.Range('E21').Value = '  +1.94%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '229.81'
$ws.Range('E22').Value = '  -3.01%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.09'
$ws.Range('E23').Value = '  -1.95%  '

# Row 24
$ws.Range('E24').Value = '  +0.23%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  -2.34%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.30'
$ws.Range('E26').Value = '  +0.27%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  -5.20%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.73'
$ws.Range('E28').Value = '  +2.51%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.44'
$ws.Range('E29').Value = '  -1.10%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.125'
$ws.Range('E30').Value = '  -4.59%  '

# Row 31
$ws.Range('E31').Value = '  -1.17%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.16'
$ws.Range('E32').Value = '  +1.82%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.70'
$ws.Range('E33').Value = '  -3.26%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0629'
$ws.Range('E34').Value = '  +1.78%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.30'
$ws.Range('E35').Value = '  -1.89%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.24'
$ws.Range('E36').Value = '  -0.63%  '

# Row 37
$ws.Range('E37').Value = '  -0.28%  '

# Row 38
$ws.Range('E38').Value = '  -3.41%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.15'
$ws.Range('E39').Value = '  -5.38%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.05'
$ws.Range('E40').Value = '  -1.81%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0987'
$ws.Range('E41').Value = '  +0.53%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.89'
$ws.Range('E42').Value = '  +3.34%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.18'
$ws.Range('E43').Value = '  -3.17%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0210'
$ws.Range('E44').Value = '  -0.57%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '16.16'
$ws.Range('E45').Value = '  +0.16%  '

# Row 46
$ws.Range('D46').Value = '1.353.18'
$ws.Range('E46').Value = '  +1.20%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.04'
$ws.Range('E47').Value = '  -4.46%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '88.08'
$ws.Range('E48').Value = '  -4.40%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.19'
$ws.Range('E49').Value = '  -4.62%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.82'
$ws.Range('E50').Value = '  +2.10%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '45.80'
$ws.Range('E51').Value = '  +4.70%  '
